# feat: add 2022-Q3 data
#
# The workbook originally has three sheets: "总计" (summary), "2022-Q2",
# "2022-Q1". This change inserts a new "2022-Q3" sheet (fund holdings
# table, same layout as the quarterly sheets) right after "总计" and
# before the former "2022-Q2" sheet, and updates the "总计" roll-up
# sheet so it lists all three quarters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it
#    inherits the exact same column layout/styles), then overwrite its
#    data cells with the Q3 figures. We track sheets by Index rather
#    than by guessing Excel's auto-generated copy name, and we re-fetch
#    the worksheet object after Move() since the old COM reference does
#    not track the sheet across a position change.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2Index = $q2.Index
$q2.Copy($null, $q2)
$q3 = $wb.Worksheets.Item($q2Index + 1)
$q3.Name = "2022-Q3"
$q3.Move($q2)
$q3 = $wb.Worksheets.Item("2022-Q3")

# Columns B (fund code) and D/E/F/G (numeric-looking figures stored as
# text in this workbook, e.g. "8.53", "94.80", "0.7413") must keep their
# text type -- force text format so Excel doesn't coerce them (and drop
# the leading zeros on fund codes) when the values are assigned below.
$q3.Range("B2:B11").NumberFormat = "@"
$q3.Range("D2:G11").NumberFormat = "@"

# Header row (row 1) is unchanged: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名

# Row 2 : 011201 财通优势行业轮动混合A
$q3.Range("B2").Value = "011201"
$q3.Range("C2").Value = "财通优势行业轮动混合A"
$q3.Range("D2").Value = "8.53"
$q3.Range("E2").Value = "91.75"
$q3.Range("F2").Value = "8.69"
$q3.Range("G2").Value = "0.7413"
$q3.Range("H2").Value = 1

# Row 3 : 010418 财通景气行业混合A
$q3.Range("B3").Value = "010418"
$q3.Range("C3").Value = "财通景气行业混合A"
$q3.Range("D3").Value = "2.72"
$q3.Range("E3").Value = "94.88"
$q3.Range("F3").Value = "9.69"
$q3.Range("G3").Value = "0.2636"
$q3.Range("H3").Value = 1

# Row 4 : 501015 财通多策略升级混合（LOF）A
$q3.Range("B4").Value = "501015"
$q3.Range("C4").Value = "财通多策略升级混合（LOF）A"
$q3.Range("D4").Value = "2.06"
$q3.Range("E4").Value = "94.80"
$q3.Range("F4").Value = "8.21"
$q3.Range("G4").Value = "0.1691"
$q3.Range("H4").Value = 5

# Row 5 : 005959 财通新视野灵活配置混合C
$q3.Range("B5").Value = "005959"
$q3.Range("C5").Value = "财通新视野灵活配置混合C"
$q3.Range("D5").Value = "1.12"
$q3.Range("E5").Value = "94.59"
$q3.Range("F5").Value = "8.99"
$q3.Range("G5").Value = "0.1007"
$q3.Range("H5").Value = 1

# Row 6 : 005851 财通新视野灵活配置混合A
$q3.Range("B6").Value = "005851"
$q3.Range("C6").Value = "财通新视野灵活配置混合A"
$q3.Range("D6").Value = "0.62"
$q3.Range("E6").Value = "94.59"
$q3.Range("F6").Value = "8.99"
$q3.Range("G6").Value = "0.0557"
$q3.Range("H6").Value = 1

# Row 7 : 501032 财通福盛多策略混合（LOF）A
$q3.Range("B7").Value = "501032"
$q3.Range("C7").Value = "财通福盛多策略混合（LOF）A"
$q3.Range("D7").Value = "0.56"
$q3.Range("E7").Value = "93.59"
$q3.Range("F7").Value = "9.14"
$q3.Range("G7").Value = "0.0512"
$q3.Range("H7").Value = 1

# Row 8 : 015271 财通多策略升级混合（LOF）C
$q3.Range("B8").Value = "015271"
$q3.Range("C8").Value = "财通多策略升级混合（LOF）C"
$q3.Range("D8").Value = "0.62"
$q3.Range("E8").Value = "94.80"
$q3.Range("F8").Value = "8.21"
$q3.Range("G8").Value = "0.0509"
$q3.Range("H8").Value = 5

# Row 9 : 011202 财通优势行业轮动混合C
$q3.Range("B9").Value = "011202"
$q3.Range("C9").Value = "财通优势行业轮动混合C"
$q3.Range("D9").Value = "0.34"
$q3.Range("E9").Value = "91.75"
$q3.Range("F9").Value = "8.69"
$q3.Range("G9").Value = "0.0295"
$q3.Range("H9").Value = 1

# Row 10 : 014628 财通福盛多策略混合（LOF）C
$q3.Range("B10").Value = "014628"
$q3.Range("C10").Value = "财通福盛多策略混合（LOF）C"
$q3.Range("D10").Value = "0.09"
$q3.Range("E10").Value = "93.59"
$q3.Range("F10").Value = "9.14"
$q3.Range("G10").Value = "0.0082"
$q3.Range("H10").Value = 1

# Row 11 : 016234 财通景气行业混合C (market value column stays a real
# number 0 here, same as it did on the "2022-Q2" sheet it was copied
# from)
$q3.Range("B11").Value = "016234"
$q3.Range("C11").Value = "财通景气行业混合C"
$q3.Range("D11").Value = "0.00"
$q3.Range("E11").Value = "94.88"
$q3.Range("F11").Value = "9.69"
$q3.Range("G11").NumberFormat = "General"
$q3.Range("G11").Value = 0
$q3.Range("H11").Value = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" roll-up sheet: it now needs a row for 2022-Q3,
#    2022-Q2 and 2022-Q1 (previously only had Q2/Q1, with the old
#    "2022-Q2" row now becoming the Q3 one).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make sure the newly used A3/A4 cells pick up the same style as A2
# (centered/bold per the "2" cell style) before writing their values.
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A2").Copy($total.Range("A4"))

# Row 4 (was row 3): 2022-Q1, unchanged figures, index bumped to 2
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 13
$total.Range("D4").Value = 1.74

# Row 3 (new): 2022-Q2, the figures that used to sit in row 2
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 10
$total.Range("D3").Value = 2.01

# Row 2: now describes 2022-Q3
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 1.47
